$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Every exact-match "int f0/1" becomes "int fa0/1", typed by placing the
#    cursor between the "f" and the "0" and inserting an "a" -- this leaves
#    the paragraph holding three runs: "int f" / "a" / "0/1".
# ---------------------------------------------------------------------------
$searchRange = $d.Content
while ($true) {
    $searchRange.Find.Execute("int f0/1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $searchRange.Find.Found) { break }
    $hitStart = $searchRange.Start

    $ip = $d.Range($hitStart + 5, $hitStart + 5)
    $ip.InsertBefore("a")

    $subA = $d.Range($hitStart + 5, $hitStart + 6)
    $subA.Font.Bold = 1
    $subA.Font.Bold = 0

    $searchRange = $d.Range($hitStart + 9, $hitStart + 9)
}

# ---------------------------------------------------------------------------
# 2) Insert the VLAN-naming block right before the "Int fa0/18" paragraph
#    (the one whose text is built from "Int fa0/" + "18").
# ---------------------------------------------------------------------------
$searchRange = $d.Content
$target = $null
while ($true) {
    $searchRange.Find.Execute("Int fa0/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $searchRange.Find.Found) { break }
    $p = $searchRange.Paragraphs(1)
    if ($p.Range.Text -eq "Int fa0/18`r") {
        $target = $p
        break
    }
    $searchRange.Collapse(0)
}

$lines1 = @(
    "Vlan 10",
    "Name Operations",
    "Vlan 20",
    "Name Parking_Lot",
    "Vlan 99",
    "Name Management",
    "Vlan 1000",
    "Name Native",
    "End",
    "Conf t"
)
$insPoint = $d.Range($target.Range.Start, $target.Range.Start)
$insPoint.InsertBefore(($lines1 -join "`r") + "`r")

# ---------------------------------------------------------------------------
# 3) Append a new "trunk native vlan 1000" sequence after the very last
#    "conf t / int f0/1(->fa0/1) / switchport mode trunk / end" block.
#    ("switchport mode trunk" also flattens-matches the earlier, two-run
#    "switchport mode" + " trunk" paragraph, so walk to the LAST hit.)
# ---------------------------------------------------------------------------
$searchRange = $d.Content
$trunkPara = $null
while ($true) {
    $searchRange.Find.Execute("switchport mode trunk", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $searchRange.Find.Found) { break }
    $trunkPara = $searchRange.Paragraphs(1)
    $searchRange.Collapse(0)
}
$endPara = $trunkPara.Next()

$lines2 = @(
    "conf t",
    "int fa0/1",
    "switchport mode dynamic desirable",
    "end",
    "conf t",
    "conf t",
    "int fa0/1",
    "switchport trunk native vlan 1000",
    "end",
    ""
)
$appendPoint = $d.Range($endPara.Range.End, $endPara.Range.End)
$appendPoint.InsertAfter(($lines2 -join "`r"))

Write-Host "Edit complete."
